# Rediseno de GUI, con CustomTkinter, y reposicion de buttons
#
# The underlying dataset (columns A-J) is a 15-row employee sample that is
# replicated throughout the sheet (rows 2-271). Column K holds the cluster
# id assigned to each employee by the clustering algorithm. The clustering
# run was redone and the cluster ids were renumbered; because each of the
# 15 distinct employee records repeats every 15 rows, the new cluster id
# only depends on the employee's position within that 15-row cycle, so the
# renumbering can be reproduced with a simple lookup table keyed by that
# cycle position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cluster id, indexed by position within the 15-row repeating cycle
# (position 0 -> row 2, position 1 -> row 3, ... position 14 -> row 16,
# then the cycle repeats starting at row 17, and so on).
$newClusterByPosition = @(0, 2, 4, 3, 4, 2, 0, 2, 1, 1, 3, 1, 3, 0, 0)
$cycleLength = $newClusterByPosition.Length

$firstDataRow = 2
$lastDataRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastDataRow -lt $firstDataRow) {
    $lastDataRow = 271
}

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $position = ($row - $firstDataRow) % $cycleLength
    $ws.Cells.Item($row, 11).Value = $newClusterByPosition[$position]
}
